$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "MFD00398B" record (row 25), shifting all
# subsequent rows up by one. This re-aligns the fieldsample_barcode
# sequence id map (MFD00399..MFD00405 -> rows 25..31).
$ws.Rows.Item(25).Delete()
